$wb = $excel.ActiveWorkbook

# --- Rename existing sheets ---
$dataSet = $wb.Worksheets.Item("Davis")
$dataSet.Name = "Data_Set"

$descSet = $wb.Worksheets.Item("definitions")
$descSet.Name = "Description_of_Data_Set"

# --- Reorder: Description_of_Data_Set, Data_Set ---
$descSet.Move($dataSet)

# NOTE: Move() swaps which physical sheet a variable handle refers to, so
# re-fetch live handles by (new) name immediately afterwards.
$dataSet = $wb.Worksheets.Item("Data_Set")
$descSet = $wb.Worksheets.Item("Description_of_Data_Set")

# --- Add the two new (empty) sheets at the end, in order so sheetId 3/4 match ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$n25 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$n25.Name = "Put n = 25 here"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$n50 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$n50.Name = "Put n = 50 here"

# --- Update sheet view / selection state ---
# Data_Set: scroll back to top and select K12
$dataSet.Activate()
$dataSet.Range("K12").Select()

# Put n = 50 here: select I24
$n50.Activate()
$n50.Range("I24").Select()

# Description_of_Data_Set: select Q15, and leave it the active tab
$descSet.Activate()
$descSet.Range("Q15").Select()

# --- Update the defined name / filter database range to point at Data_Set ---
$wb.Names("_xlnm._FilterDatabase").RefersTo = "=Data_Set!`$B`$1:`$F`$200"

